$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. It belongs right after the
# existing row 199, so insert a fresh row at 200 (pushing every following
# row down by one, including the final row which becomes row 236).
$ws.Rows(200).Insert()

# The new record shares the market/category/variety/quality/unit/origin
# attributes with what is now row 201 (the old row 200) - copy that row
# down into the freshly inserted row 200 so every column is populated
# (and formatting, e.g. the date style on column D, comes along too).
$ws.Range("A201:R201").Copy()
$ws.Range("A200").PasteSpecial()

# Now overwrite just the figures that differ for this new observation.
$ws.Range("D200").Value = 44776
$ws.Range("J200").Value = 580
$ws.Range("K200").Value = 17000
$ws.Range("M200").Value = 17448
$ws.Range("P200").Value = 698
